$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the K-column formulas in each of the 7 data blocks so that they reference
# the block's own sample-size cell (G2, G14, G29, G47, G68, G84, G95) instead of
# the hardcoded literal 10500.
$blocks = @(
    @{ Start = 4;  End = 9;   Ref = "G2" },
    @{ Start = 16; End = 24;  Ref = "G14" },
    @{ Start = 31; End = 41;  Ref = "G29" },
    @{ Start = 49; End = 62;  Ref = "G47" },
    @{ Start = 70; End = 78;  Ref = "G68" },
    @{ Start = 86; End = 89;  Ref = "G84" },
    @{ Start = 97; End = 107; Ref = "G95" }
)

foreach ($b in $blocks) {
    for ($r = $b.Start; $r -le $b.End; $r++) {
        $ws.Range("K$r").Formula = "=LOG($" + $b.Ref.Substring(0,1) + "$" + $b.Ref.Substring(1) + ")"
    }
}

# Extend the AVERAGE/MIN/MAX summary ranges for the blocks whose ranges
# previously stopped short of the last data row (31-41, 49-62, 97-107).
$ws.Range("M42").Formula = "=AVERAGE(M31:M41)"
$ws.Range("M43").Formula = "=MIN(M31:M41)"
$ws.Range("M44").Formula = "=MAX(M31:M41)"
$ws.Range("O42").Formula = "=AVERAGE(O31:O41)"
$ws.Range("O43").Formula = "=MIN(O31:O41)"
$ws.Range("O44").Formula = "=MAX(O31:O41)"

$ws.Range("M63").Formula = "=AVERAGE(M49:M62)"
$ws.Range("M64").Formula = "=MIN(M49:M62)"
$ws.Range("M65").Formula = "=MAX(M49:M62)"
$ws.Range("O63").Formula = "=AVERAGE(O49:O62)"
$ws.Range("O64").Formula = "=MIN(O49:O62)"
$ws.Range("O65").Formula = "=MAX(O49:O62)"

$ws.Range("M108").Formula = "=AVERAGE(M97:M107)"
$ws.Range("M109").Formula = "=MIN(M97:M107)"
$ws.Range("M110").Formula = "=MAX(M97:M107)"
$ws.Range("O108").Formula = "=AVERAGE(O97:O107)"
$ws.Range("O109").Formula = "=MIN(O97:O107)"
$ws.Range("O110").Formula = "=MAX(O97:O107)"

$excel.Calculate()
